# ==========================================================================
# SOLUCION DE CONFLICTOS EXCEL Y GESTION USUARIOS INICIO/INACTIVOS
# - Append a brand-new row (row 8) for the Jose Gonzalez admin user
# - Replace the Miguel Garcia admin record (row 2) with the same Jose
#   Gonzalez data (duplicate/merge of the conflicting admin accounts)
# - Mark the Julia Ruiz (row 5) and Carolina Castro (row 7) accounts active
# ==========================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new user record for Jose Gonzalez (typed in first) ---
[void]$ws.Range("A2:L2").Copy()
[void]$ws.Range("A8:L8").PasteSpecial(-4122)
$ws.Range("A8").Value2 = "jose.gonzalez&pineed"
$ws.Range("B8").Value2 = 3256451010312
$ws.Range("C8").Value2 = "José"
$ws.Range("D8").Value2 = "González"
$ws.Range("E8").Value2 = "ADMINISTRADOR"
$ws.Range("F8").Value2 = "Masculino"
$ws.Range("G8").Value2 = 3256451010312
$ws.Range("H8").Value2 = 31168
$ws.Range("I8").Value2 = 59596565
$ws.Range("J8").Value2 = "jose.gonzalez@gmail.com"
$ws.Range("J8").Style = "Normal"
$ws.Range("K8").Value2 = "ACTIVO"
$ws.Range("L8").Value2 = $true
$ws.Rows.Item(8).RowHeight = 30

# --- Row 2: update Nombre/Apellido/DPI/Telefono/Correo for the admin user ---
$ws.Range("C2").Value2 = "José"
$ws.Range("D2").Value2 = "González"
$ws.Range("G2").Value2 = 3256451010312
$ws.Range("I2").Value2 = 59596565
$ws.Range("J2").Value2 = "jose.gonzalez@gmail.com"
$ws.Range("J2").Style = "Normal"
$ws.Rows.Item(2).RowHeight = 30

# --- Rows 5 & 7: flip their "activo" flag (column L) on ---
$ws.Range("L5").Value2 = $true
$ws.Range("L7").Value2 = $true

[void]$ws.Range("M4").Select()

Write-Host "edit complete"
